$d = $word.ActiveDocument

# The three "signature" paragraphs (name + title separated by a line break)
# are being replaced by a 3-row x 2-column table (name cell + title cell per
# inspector), preceded by a blank paragraph. Locate the three paragraphs by
# their distinctive text.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Sanjay Kumar Singh*") {
        $startPara = $i
    }
    if ($t -like "Vivek Kumar*") {
        $endPara = $i
        break
    }
}

$startPos = $d.Paragraphs.Item($startPara).Range.Start
$endPos = $d.Paragraphs.Item($endPara).Range.End

# Remove the three paragraphs (including their paragraph marks) entirely.
$killRange = $d.Range($startPos, $endPos)
$killRange.Delete()

# Build the replacement: an empty paragraph followed by the new table,
# expressed as raw WordprocessingML so the exact formatting/markup (e.g.
# literal percentage widths) is reproduced faithfully.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-NameCell([string]$name, [string]$jc, [bool]$isTopCell) {
    if ($isTopCell) {
        $tcMar = '<w:tcMar><w:top w:type="dxa" w:w="100"/><w:left w:type="dxa" w:w="100"/><w:bottom w:type="dxa" w:w="50"/><w:right w:type="dxa" w:w="100"/></w:tcMar>'
    } else {
        $tcMar = '<w:tcMar><w:top w:type="dxa" w:w="50"/><w:left w:type="dxa" w:w="100"/><w:bottom w:type="dxa" w:w="100"/><w:right w:type="dxa" w:w="100"/></w:tcMar>'
    }
    $tcBorders = '<w:tcBorders><w:top w:val="none"/><w:left w:val="none"/><w:bottom w:val="none"/><w:right w:val="none"/></w:tcBorders>'
    return '<w:tc><w:tcPr><w:tcW w:type="pct" w:w="50%"/>' + $tcBorders + $tcMar + '</w:tcPr>' +
        '<w:p><w:pPr><w:spacing w:after="50"/><w:jc w:val="' + $jc + '"/></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">' + $name + '</w:t></w:r></w:p></w:tc>'
}

function New-TitleCell([string]$title, [string]$jc, [bool]$isTopCell) {
    if ($isTopCell) {
        $tcMar = '<w:tcMar><w:top w:type="dxa" w:w="100"/><w:left w:type="dxa" w:w="100"/><w:bottom w:type="dxa" w:w="50"/><w:right w:type="dxa" w:w="100"/></w:tcMar>'
    } else {
        $tcMar = '<w:tcMar><w:top w:type="dxa" w:w="50"/><w:left w:type="dxa" w:w="100"/><w:bottom w:type="dxa" w:w="100"/><w:right w:type="dxa" w:w="100"/></w:tcMar>'
    }
    $tcBorders = '<w:tcBorders><w:top w:val="none"/><w:left w:val="none"/><w:bottom w:val="none"/><w:right w:val="none"/></w:tcBorders>'
    return '<w:tc><w:tcPr><w:tcW w:type="pct" w:w="50%"/>' + $tcBorders + $tcMar + '</w:tcPr>' +
        '<w:p><w:pPr><w:spacing w:after="50"/><w:jc w:val="' + $jc + '"/></w:pPr>' +
        '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">' + $title + '</w:t></w:r></w:p></w:tc>'
}

function New-InspectorRow([string]$name, [string]$title, [string]$jc) {
    return '<w:tr>' + (New-NameCell $name $jc $true) + (New-TitleCell $title $jc $false) + '</w:tr>'
}

$emptyPara = '<w:p ' + $wNs + '><w:pPr><w:spacing w:after="300"/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p>'

$tblPr = '<w:tblPr><w:tblW w:type="pct" w:w="100%"/><w:tblBorders>' +
    '<w:top w:val="single" w:color="auto" w:sz="4"/>' +
    '<w:left w:val="single" w:color="auto" w:sz="4"/>' +
    '<w:bottom w:val="single" w:color="auto" w:sz="4"/>' +
    '<w:right w:val="single" w:color="auto" w:sz="4"/>' +
    '<w:insideH w:val="single" w:color="auto" w:sz="4"/>' +
    '<w:insideV w:val="single" w:color="auto" w:sz="4"/>' +
    '</w:tblBorders></w:tblPr>' +
    '<w:tblGrid><w:gridCol w:w="100"/><w:gridCol w:w="100"/></w:tblGrid>'

$rows = (New-InspectorRow "Sanjay Kumar Singh" "CMI/YTSK" "right") +
    (New-InspectorRow "Lovey Gandhi" "CMI/G." "center") +
    (New-InspectorRow "Vivek Kumar" "CMI/Ctg/VIP" "left")

$table = '<w:tbl ' + $wNs + '>' + $tblPr + $rows + '</w:tbl>'

$insertRange = $d.Range($startPos, $startPos)
$insertRange.InsertXML($emptyPara + $table) | Out-Null
